$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 356.075
$ws.Range("I80").Value = 294.3
$ws.Range("J80").Value = 541.4
$ws.Range("K80").Value = 882.9000000000001
$ws.Range("L80").Value = 1624.2
$ws.Range("M80").Value = 115.0999999999999
$ws.Range("N80").Value = -3620.2
$ws.Range("H83").Value = 356.075
$ws.Range("I83").Value = 294.3
$ws.Range("J83").Value = 541.4
$ws.Range("K83").Value = 2648.7
$ws.Range("L83").Value = 4872.599999999999
$ws.Range("M83").Value = 2343.3
$ws.Range("N83").Value = -14856.6
$ws.Range("H129").Value = 1893.7646
$ws.Range("J129").Value = 2521.7
$ws.Range("L129").Value = 7565.099999999999
$ws.Range("N129").Value = -17565.1
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2453054.5
$ws.Range("I2").Value = 3638.2
$ws.Range("J2").Value = 4202637.5
$ws.Range("K2").Value = 3638.2
$ws.Range("L2").Value = 4202637.5
$ws.Range("M2").Value = -3525.2
$ws.Range("N2").Value = -4202863.5
$ws.Range("H45").Value = 902.9091
$ws.Range("I45").Value = 846.5
$ws.Range("J45").Value = 1053.3334
$ws.Range("K45").Value = 846.5
$ws.Range("L45").Value = 1053.3334
$ws.Range("M45").Value = -469.5
$ws.Range("N45").Value = -1807.3334
$ws.Range("H61").Value = 2053.2144
$ws.Range("I61").Value = 1428.85
$ws.Range("J61").Value = 3614.125
$ws.Range("K61").Value = 1428.85
$ws.Range("L61").Value = 3614.125
$ws.Range("M61").Value = -1216.85
$ws.Range("N61").Value = -4038.125
$ws.Range("H74").Value = 3270.5293
$ws.Range("I74").Value = 667.75
$ws.Range("J74").Value = 6439.1304
$ws.Range("K74").Value = 667.75
$ws.Range("L74").Value = 6439.1304
$ws.Range("M74").Value = 206.25
$ws.Range("N74").Value = -8187.1304
$ws.Range("H77").Value = 3270.5293
$ws.Range("I77").Value = 667.75
$ws.Range("J77").Value = 6439.1304
$ws.Range("K77").Value = 3338.75
$ws.Range("L77").Value = 32195.652
$ws.Range("M77").Value = 1029.25
$ws.Range("N77").Value = -40931.652
$ws.Range("H88").Value = 3225
$ws.Range("I88").Value = 3225
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 3225
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -2819
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 3225
$ws.Range("I91").Value = 3225
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 3225
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -1821
$ws.Range("N91").ClearContents()
$ws.Range("H96").Value = 33997.668
$ws.Range("J96").Value = 33997.668
$ws.Range("L96").Value = 33997.668
$ws.Range("N96").Value = -39489.668
$ws.Range("H102").Value = 1479.2354
$ws.Range("I102").Value = 1512
$ws.Range("K102").Value = 1512
$ws.Range("M102").Value = 110
$ws.Range("H110").Value = 1108.2858
$ws.Range("I110").Value = 971.1667
$ws.Range("J110").Value = 1931
$ws.Range("K110").Value = 971.1667
$ws.Range("L110").Value = 1931
$ws.Range("M110").Value = 1073.8333
$ws.Range("N110").Value = -6021
$ws.Range("H116").Value = 2453054.5
$ws.Range("I116").Value = 3638.2
$ws.Range("J116").Value = 4202637.5
$ws.Range("K116").Value = 3638.2
$ws.Range("L116").Value = 4202637.5
$ws.Range("M116").Value = -1344.2
$ws.Range("N116").Value = -4207225.5
$ws.Range("H132").Value = 2734.7188
$ws.Range("I132").Value = 2672.0476
$ws.Range("J132").Value = 2854.3635
$ws.Range("K132").Value = 8016.1428
$ws.Range("L132").Value = 8563.0905
$ws.Range("M132").Value = -5486.1428
$ws.Range("N132").Value = -13623.0905
$ws.Range("H136").Value = 2053.2144
$ws.Range("I136").Value = 1428.85
$ws.Range("J136").Value = 3614.125
$ws.Range("K136").Value = 4286.549999999999
$ws.Range("L136").Value = 10842.375
$ws.Range("M136").Value = -1736.549999999999
$ws.Range("N136").Value = -15942.375
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2453054.5
$ws.Range("I3").Value = 3638.2
$ws.Range("J3").Value = 4202637.5
$ws.Range("K3").Value = 3638.2
$ws.Range("L3").Value = 4202637.5
$ws.Range("M3").Value = -3524.2
$ws.Range("N3").Value = -4202865.5
$ws.Range("H50").Value = 29326.666
$ws.Range("J50").Value = 29326.666
$ws.Range("L50").Value = 29326.666
$ws.Range("N50").Value = -30474.666
$ws.Range("H86").Value = 2385.6
$ws.Range("I86").Value = 2282.6155
$ws.Range("K86").Value = 2282.6155
$ws.Range("M86").Value = -1159.6155
$ws.Range("H89").Value = 2385.6
$ws.Range("I89").Value = 2282.6155
$ws.Range("K89").Value = 11413.0775
$ws.Range("M89").Value = -5797.077499999999
$ws.Range("H107").Value = 2150
$ws.Range("I107").Value = 2150
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2150
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -230
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 59437.633
$ws.Range("I134").Value = 70213.44
$ws.Range("K134").Value = 210640.32
$ws.Range("M134").Value = -208105.32
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 2036
$ws.Range("I7").Value = 2532.5
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 2532.5
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = -2419.5
$ws.Range("N7").Value = -276
$ws.Range("H22").Value = 469.06384
$ws.Range("I22").Value = 423.88095
$ws.Range("J22").Value = 848.6
$ws.Range("K22").Value = 423.88095
$ws.Range("L22").Value = 848.6
$ws.Range("M22").Value = -73.88094999999998
$ws.Range("N22").Value = -1548.6
$ws.Range("H28").Value = 15571.5
$ws.Range("J28").Value = 15571.5
$ws.Range("L28").Value = 15571.5
$ws.Range("N28").Value = -16061.5
$ws.Range("H121").Value = 31000
$ws.Range("J121").Value = 31000
$ws.Range("L121").Value = 31000
$ws.Range("N121").Value = -33620
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 2666.6667
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 2666.6667
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 8000.000100000001
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -12118.0001
$ws.Range("H105").Value = 908000000
$ws.Range("J105").Value = 908000000
$ws.Range("L105").Value = 2724000000
$ws.Range("N105").Value = -2724005242
$ws.Range("H131").Value = 1589851.6
$ws.Range("I131").Value = 4636.5
$ws.Range("J131").Value = 2042770.1
$ws.Range("K131").Value = 13909.5
$ws.Range("L131").Value = 6128310.300000001
$ws.Range("M131").Value = -8869.5
$ws.Range("N131").Value = -6138390.300000001
$ws.Range("H133").Value = 5052.4053
$ws.Range("I133").Value = 2171.1875
$ws.Range("J133").Value = 7247.619
$ws.Range("K133").Value = 6513.5625
$ws.Range("L133").Value = 21742.857
$ws.Range("M133").Value = -1453.5625
$ws.Range("N133").Value = -31862.857
$ws.Range("H139").Value = 950.9091
$ws.Range("I139").Value = 866
$ws.Range("J139").Value = 1800
$ws.Range("K139").Value = 2598
$ws.Range("L139").Value = 5400
$ws.Range("M139").Value = 2542
$ws.Range("N139").Value = -15680
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1308
$ws.Range("I102").Value = 935.1667
$ws.Range("K102").Value = 935.1667
$ws.Range("M102").Value = 686.8333
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1971.9131
$ws.Range("I40").Value = 1913.8889
$ws.Range("K40").Value = 1913.8889
$ws.Range("M40").Value = -1777.8889
$ws.Range("H61").Value = 698.4666999999999
$ws.Range("I61").Value = 605.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 605.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -403.5
$ws.Range("N61").Value = -2404
$ws.Range("H95").Value = 11555.429
$ws.Range("J95").Value = 11555.429
$ws.Range("L95").Value = 11555.429
$ws.Range("N95").Value = -17047.429
$ws.Range("H113").Value = 698.4666999999999
$ws.Range("I113").Value = 605.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 605.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 1564.5
$ws.Range("N113").Value = -6340
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 42500
$ws.Range("J92").Value = 42500
$ws.Range("L92").Value = 42500
$ws.Range("N92").Value = -47492
$ws.Range("H132").Value = 8611.444
$ws.Range("I132").Value = 14550
$ws.Range("J132").Value = 3860.6
$ws.Range("K132").Value = 43650
$ws.Range("L132").Value = 11581.8
$ws.Range("M132").Value = -41120
$ws.Range("N132").Value = -16641.8
